{"js": "const body = context.document.body;\n\n// Each entry: exact full-paragraph text to match -> replacement text.\n// Matching is scoped per-paragraph (on the paragraph's own full text) so we\n// only touch the intended bullet/heading/meta paragraphs and never a\n// substring occurring inside a longer sentence elsewhere in the document.\nconst replacements = [\n  [\"Play Cafelito Slot for Free - Review 2021\", \"Play Cafelito Free- A Slot Game Review\"],\n  [\"4806 ways to win\", \"Simple and straightforward gameplay\"],\n  [\"Affordable minimum bet at \\u20AC0.40\", \"Highly detailed and colorful graphics\"],\n  [\"Detailed and colorful graphics\", \"Special symbols and bonuses increase winnings\"],\n  [\"Several bonuses and jackpots available\", \"Three variable jackpots to strive for\"],\n  [\"Only three jackpots available\", \"Limited number of free spins\"],\n  [\"RTP value could be higher\", \"Jackpots can be difficult to win\"],\n  [\n    \"Read our review of Cafelito slot game and play for free. Discover the graphics, gameplay, bonuses, RTP value, and jackpots. Start playing now.\",\n    \"Read our review of Cafelito, a slot game with simple gameplay and three variable jackpots. Play for free now!\",\n  ],\n];\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  for (const [find, replace] of replacements) {\n    if (paragraph.text === find) {\n      const found = paragraph.search(find, { matchCase: true, matchWholeWord: false });\n      found.load(\"items\");\n      await context.sync();\n\n      for (const range of found.items) {\n        range.insertText(replace, \"Replace\");\n      }\n      await context.sync();\n      break;\n    }\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Exact full-paragraph text (Range.Text includes the trailing paragraph\n# mark, hence the `r) -> replacement text. Matching a paragraph's whole\n# text (rather than a document-wide Find) keeps this scoped to the\n# intended bullet/heading/meta paragraphs and avoids touching the same\n# substring when it appears inside a longer sentence elsewhere.\n$replacements = [ordered]@{\n    \"Play Cafelito Slot for Free - Review 2021`r\" = \"Play Cafelito Free- A Slot Game Review\"\n    \"4806 ways to win`r\" = \"Simple and straightforward gameplay\"\n    \"Affordable minimum bet at \u20ac0.40`r\" = \"Highly detailed and colorful graphics\"\n    \"Detailed and colorful graphics`r\" = \"Special symbols and bonuses increase winnings\"\n    \"Several bonuses and jackpots available`r\" = \"Three variable jackpots to strive for\"\n    \"Only three jackpots available`r\" = \"Limited number of free spins\"\n    \"RTP value could be higher`r\" = \"Jackpots can be difficult to win\"\n    \"Read our review of Cafelito slot game and play for free. Discover the graphics, gameplay, bonuses, RTP value, and jackpots. Start playing now.`r\" = \"Read our review of Cafelito, a slot game with simple gameplay and three variable jackpots. Play for free now!\"\n}\n\nforeach ($p in $d.Paragraphs) {\n    $paraText = $p.Range.Text\n    if ($replacements.Contains($paraText)) {\n        $newText = $replacements[$paraText]\n        $findText = $paraText.Substring(0, $paraText.Length - 1)\n\n        $rng = $p.Range\n        $rng.Find.ClearFormatting()\n        $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    }\n}\n"}
